$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in the title cell (A1): 08:16 -> 08:46
$ws.Range("A1").Value = "Datos actualizados a 22 de Marzo de 2020 a las 08:46"

# Hungria moves up the ranking (new data) to sit right after Uruguay (row 71),
# pushing Letonia, Costa Rica, Republica Dominicana and Lituania down by one row.
# Insert a fresh row at 72 for Hungria with its updated stats.
$ws.Rows("72:72").Insert()
$ws.Range("A72").Value = "Hungria"
$ws.Range("B72").Value = 131
$ws.Range("C72").Value = 28
$ws.Range("D72").Value = 7
$ws.Range("E72").Value = 120
$ws.Range("F72").Value = 6
$ws.Range("G72").Value = 0
$ws.Range("H72").Value = 4

# Remove the old Hungria row, which (after the insert above) now sits at row 77.
$ws.Rows("77:77").Delete()
